$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 186-189: replace the "ALL patients" description
# text with the new "patients with a valid sex, DOB..." description. The
# B/C column text (Location FIPS / Census Tract|Block Group Geocode (year))
# is unchanged - only the D column descriptions and the row formatting
# change (row height 51 -> 68; B/C lose their explicit font styling and D
# keeps the standard wrapped-text look used elsewhere in the sheet).

$ws.Range("D186").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have have a valid, 11-digit tract level geocode for the 2010 census year associated with their current location (location_id in the person table)"
$ws.Range("D187").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have a valid, 11-digit tract level geocode for the 2020 census year associated with their current location (location_id in the person table)"
$ws.Range("D188").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have a valid, 12-digit block group level geocode for the 2010 census year associated with their current location (location_id in the person table)"
$ws.Range("D189").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have a valid, 12-digit block group level geocode for the 2020 census year associated with their current location (location_id in the person table)"

foreach ($r in 186..189) {
    $ws.Range("B$r").Style = "Normal"
    $ws.Range("C$r").Style = "Normal"
    $ws.Range("D$r").Style = "Normal"
    $ws.Range("D$r").WrapText = $true
    $ws.Rows.Item($r).RowHeight = 68
}

# --- Add four new rows (190-193) for the "Location History" geocode checks.
# Clone the formatting from the now-updated row 186 (same A/B/C/D styling)
# so the new rows pick up the identical cell styles (A keeps the
# "Expected Concepts Present" shading, D keeps the wrapped-text style).

$ws.Range("A186:D186").Copy()
$ws.Range("A190:D193").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

foreach ($r in 190..193) {
    $ws.Rows.Item($r).RowHeight = 68
}

$ws.Range("A190").Value = "Expected Concepts Present"
$ws.Range("B190").Value = "Location History"
$ws.Range("C190").Value = "2+ Census Tract Location History Geocodes (2010)"
$ws.Range("D190").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have a valid, 11-digit tract level geocode for the 2010 census associated with MORE THAN ONE location in the location_history table"

$ws.Range("A191").Value = "Expected Concepts Present"
$ws.Range("B191").Value = "Location History"
$ws.Range("C191").Value = "2+ Census Tract Location History Geocodes (2020)"
$ws.Range("D191").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have a valid, 11-digit tract level geocode for the 2020 census associated with MORE THAN ONE location in the location_history table"

$ws.Range("A192").Value = "Expected Concepts Present"
$ws.Range("B192").Value = "Location History"
$ws.Range("C192").Value = "2+ Census Block Group Location History Geocodes (2010)"
$ws.Range("D192").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have a valid, 12-digit block group level geocode for the 2010 census associated with MORE THAN ONE location in the location_history table"

$ws.Range("A193").Value = "Expected Concepts Present"
$ws.Range("B193").Value = "Location History"
$ws.Range("C193").Value = "2+ Census Block Group Location History Geocodes (2020)"
$ws.Range("D193").Value = "counts the proportion of patients with a valid sex, DOB, and at least 1 diagnosis associated with a FTF visit that also have a valid, 12-digit block group level geocode for the 2020 census associated with MORE THAN ONE location in the location_history table"

# --- Update the view so the active selection matches the end of the new data.
$ws.Range("D193").Select()
